$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (previously 2021-03-29 / 44284) -> becomes 2021-04-05 / 44291 data
$ws.Range("D2").Value = 44291
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 11000
$ws.Range("P2").Value = 550

# Row 4 (previously 2021-03-22 / 44277) -> becomes 2021-03-29 / 44284 data
$ws.Range("D4").Value = 44284
$ws.Range("J4").Value = 100
$ws.Range("K4").Value = 10000
$ws.Range("L4").Value = 10000
$ws.Range("M4").Value = 10000
$ws.Range("P4").Value = 500

# Row 5 (previously 2021-04-05 / 44291) -> becomes 2021-03-22 / 44277 data
$ws.Range("D5").Value = 44277
$ws.Range("J5").Value = 150
$ws.Range("K5").Value = 11000
$ws.Range("L5").Value = 11000
$ws.Range("M5").Value = 11000
$ws.Range("P5").Value = 550
